$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A values (rows 2-11)
$aValues = @(1376, 5254, 1568, 3277, 1091, 2503, 1145, 3265, 2067, 1014)
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $aValues[$i]
    $cell.NumberFormat = "#,##0"
}

# Column C values (rows 2-11)
# Most are plain integer values with "#,##0" format.
# C4 and C10 are formulas "=0.946" and "=0.717".
# C10 and C11 use "#,##0.000" format.
$ws.Cells.Item(2, 3).Value = 2693
$ws.Cells.Item(2, 3).NumberFormat = "#,##0"

$ws.Cells.Item(3, 3).Value = 1358
$ws.Cells.Item(3, 3).NumberFormat = "#,##0"

$ws.Cells.Item(4, 3).Formula = "=0.946"

$ws.Cells.Item(5, 3).Value = 6632
$ws.Cells.Item(5, 3).NumberFormat = "#,##0"

$ws.Cells.Item(6, 3).Value = 3059
$ws.Cells.Item(6, 3).NumberFormat = "#,##0"

$ws.Cells.Item(7, 3).Value = 1201
$ws.Cells.Item(7, 3).NumberFormat = "#,##0"

$ws.Cells.Item(8, 3).Value = 1325
$ws.Cells.Item(8, 3).NumberFormat = "#,##0"

$ws.Cells.Item(9, 3).Value = 1274
$ws.Cells.Item(9, 3).NumberFormat = "#,##0"

$ws.Cells.Item(10, 3).Formula = "=0.717"
$ws.Cells.Item(10, 3).NumberFormat = "#,##0.000"

$ws.Cells.Item(11, 3).Value = 0.804
$ws.Cells.Item(11, 3).NumberFormat = "#,##0.000"

# Row 13: MEDIAN formulas
$ws.Cells.Item(13, 1).Formula = "=MEDIAN(A2:A11)"
$ws.Cells.Item(13, 1).NumberFormat = "#,##0"

$ws.Cells.Item(13, 3).Formula = "=MEDIAN(C2:C11)"
$ws.Cells.Item(13, 3).NumberFormat = "#,##0"

# Update selection to match the target state
$ws.Range("I7").Select()
